# Auto-generated script to update F-column (想去人数 / attendance count) values
# across the four worksheets, matching the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9091
$ws.Range("F4").Value = 6642
$ws.Range("F5").Value = 181
$ws.Range("F6").Value = 2160
$ws.Range("F7").Value = 611
$ws.Range("F8").Value = 81
$ws.Range("F13").Value = 15
$ws.Range("F14").Value = 88
$ws.Range("F15").Value = 31
$ws.Range("F16").Value = 9068
$ws.Range("F19").Value = 205
$ws.Range("F20").Value = 119
$ws.Range("F21").Value = 1857
$ws.Range("F25").Value = 107
$ws.Range("F27").Value = 205
$ws.Range("F29").Value = 23
$ws.Range("F30").Value = 82
$ws.Range("F31").Value = 571
$ws.Range("F32").Value = 38
$ws.Range("F33").Value = 50
$ws.Range("F34").Value = 555
$ws.Range("F35").Value = 2418
$ws.Range("F36").Value = 887
$ws.Range("F37").Value = 562
$ws.Range("F41").Value = 314
$ws.Range("F42").Value = 188
$ws.Range("F45").Value = 33
$ws.Range("F47").Value = 29
$ws.Range("F48").Value = 4008
$ws.Range("F49").Value = 16

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 443
$ws.Range("F14").Value = 22

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2345
$ws.Range("F4").Value = 349
$ws.Range("F5").Value = 28

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2345
$ws.Range("F3").Value = 9090
$ws.Range("F6").Value = 6642
$ws.Range("F8").Value = 611
$ws.Range("F9").Value = 81
$ws.Range("F14").Value = 28
$ws.Range("F15").Value = 15
$ws.Range("F16").Value = 88
$ws.Range("F17").Value = 9069
$ws.Range("F20").Value = 205
$ws.Range("F21").Value = 119
$ws.Range("F22").Value = 1857
$ws.Range("F24").Value = 107
$ws.Range("F26").Value = 205
$ws.Range("F27").Value = 23
$ws.Range("F29").Value = 571
$ws.Range("F30").Value = 38
$ws.Range("F31").Value = 50
$ws.Range("F32").Value = 556
$ws.Range("F33").Value = 887
$ws.Range("F34").Value = 22
$ws.Range("F36").Value = 562
$ws.Range("F37").Value = 314
$ws.Range("F39").Value = 188
$ws.Range("F42").Value = 33
$ws.Range("F44").Value = 29
$ws.Range("F45").Value = 4008
$ws.Range("F48").Value = 16
